# Clear repeated hierarchical label cells (columns A-E) in rows 9-70 of the
# "table" worksheet. The source table repeats the category / sub-category
# labels (A, B, C, D, E) down the rows for readability, but the canonical
# OOXML export only keeps the label on the first row of each run and blanks
# it out on the following rows where it repeats the cell directly above
# (matching the JSON-record "don't repeat a parent key" flattening used to
# build this sheet). This script reproduces that: every (row, column) pair
# below was a label cell whose text equalled the cell directly above it in
# the original workbook, so it gets cleared to an empty string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = ""
$ws.Range("A10").Value = ""
$ws.Range("A11").Value = ""
$ws.Range("A12:B12").Value = ""
$ws.Range("A13:C13").Value = ""
$ws.Range("A14:C14").Value = ""
$ws.Range("A15:C15").Value = ""
$ws.Range("A16:C16").Value = ""
$ws.Range("A17:C17").Value = ""
$ws.Range("A18:B18").Value = ""
$ws.Range("A19:C19").Value = ""
$ws.Range("A20:C20").Value = ""
$ws.Range("A22:B22").Value = ""
$ws.Range("A23:C23").Value = ""
$ws.Range("A24:B24").Value = ""
$ws.Range("A25:C25").Value = ""
$ws.Range("A26:B26").Value = ""
$ws.Range("A27:B27").Value = ""
$ws.Range("A28").Value = ""
$ws.Range("A29").Value = ""
$ws.Range("A30:B30").Value = ""
$ws.Range("A31:B31").Value = ""
$ws.Range("A32:B32").Value = ""
$ws.Range("A33:B33").Value = ""
$ws.Range("A34:C34").Value = ""
$ws.Range("A35:D35").Value = ""
$ws.Range("A36:D36").Value = ""
$ws.Range("A37:D37").Value = ""
$ws.Range("A38:D38").Value = ""
$ws.Range("A39:D39").Value = ""
$ws.Range("A40").Value = ""
$ws.Range("A41:B41").Value = ""
$ws.Range("A42:B42").Value = ""
$ws.Range("A43:C43").Value = ""
$ws.Range("A44").Value = ""
$ws.Range("A45:B45").Value = ""
$ws.Range("A46:B46").Value = ""
$ws.Range("A47:C47").Value = ""
$ws.Range("A48").Value = ""
$ws.Range("A49:B49").Value = ""
$ws.Range("A50:B50").Value = ""
$ws.Range("A51:C51").Value = ""
$ws.Range("A52").Value = ""
$ws.Range("A53:B53").Value = ""
$ws.Range("A54:B54").Value = ""
$ws.Range("A55").Value = ""
$ws.Range("A56").Value = ""
$ws.Range("A57").Value = ""
$ws.Range("A58").Value = ""
$ws.Range("A59").Value = ""
$ws.Range("A60:C60").Value = ""
$ws.Range("A61:C61").Value = ""
$ws.Range("A62:C62").Value = ""
$ws.Range("A63:B63").Value = ""
$ws.Range("A64:C64").Value = ""
$ws.Range("A66").Value = ""
$ws.Range("A67:B67").Value = ""
$ws.Range("A68:B68").Value = ""
$ws.Range("A70").Value = ""
